# Generate Report for Handback
# Refreshes the localization-status report: marks handback as in sync,
# updates the latest handback timestamps, and clears the now-stale
# "version mismatch" error details for the zh-cn and de-de targets.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# --- Overview sheet -------------------------------------------------
# Status columns (zh-cn / de-de) for both file rows move from
# "Ready for handoff" to "Handed back: in sync with en-US".
$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"

$wsOverview.Columns.Item(5).ColumnWidth = 29.1666666666667
$wsOverview.Columns.Item(6).ColumnWidth = 29.1666666666667

# --- zh-cn sheet ------------------------------------------------------
# Status column for both rows also flips to "Handed back: in sync with en-US".
$wsZhCn.Range("C2").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("C3").Value = "Handed back: in sync with en-US"

# Latest handback datetime refreshed.
$wsZhCn.Range("K2").Value = "2016-08-13 08:42:11"
$wsZhCn.Range("K3").Value = "2016-08-13 08:42:11"

# Error Detail no longer applicable - handback is now current.
$wsZhCn.Range("P2").Value = ""
$wsZhCn.Range("P3").Value = ""

$wsZhCn.Columns.Item(3).ColumnWidth = 29.1666666666667
$wsZhCn.Columns.Item(16).ColumnWidth = 12.8333333333333

# --- de-de sheet --------------------------------------------------------
$wsDeDe.Range("C2").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("C3").Value = "Handed back: in sync with en-US"

$wsDeDe.Range("K2").Value = "2016-08-13 08:42:20"
$wsDeDe.Range("K3").Value = "2016-08-13 08:42:20"

$wsDeDe.Range("P2").Value = ""
$wsDeDe.Range("P3").Value = ""

$wsDeDe.Columns.Item(3).ColumnWidth = 29.1666666666667
$wsDeDe.Columns.Item(16).ColumnWidth = 12.8333333333333
